$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# --- 1) Refresh the cached "datetimeFigureOut" date text (master + every layout) ---
$newDate = "7/1/2022"
$ppPlaceholderDate = 16

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Update-DateShape $m.Shapes
for ($L = 1; $L -le $m.CustomLayouts.Count; $L++) {
    $layout = $m.CustomLayouts.Item($L)
    Update-DateShape $layout.Shapes
}

# --- 2) Shrink the Slide Master default text style font sizes ---
$ppTitleStyle = 2
$ppBodyStyle = 3

$titleStyle = $m.TextStyles.Item($ppTitleStyle)
$titleLvl1 = $titleStyle.Levels(1)
$titleLvl1.Font.Size = 28

$bodyStyle = $m.TextStyles.Item($ppBodyStyle)
$bodyStyle.Levels(1).Font.Size = 18
$bodyStyle.Levels(2).Font.Size = 18
$bodyStyle.Levels(3).Font.Size = 14
$bodyStyle.Levels(4).Font.Size = 12
$bodyStyle.Levels(5).Font.Size = 12
